# Insert a new data row at row 171, pushing the existing rows 171-279 down
# to 172-280 (all their data moves as-is), then populate the new row 171
# with its own values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(171).Insert()

$ws.Range("A171").Value = 3
$ws.Range("B171").Value = "Femacal de La Calera"
$ws.Range("C171").Value = "Coquimbo"
$ws.Range("D171").Value = 44582
$ws.Range("E171").Value = 5
$ws.Range("F171").Value = 100112043
$ws.Range("G171").Value = "Pepino ensalada"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 85
$ws.Range("K171").Value = 11000
$ws.Range("L171").Value = 12000
$ws.Range("M171").Value = 11471
$ws.Range("N171").Value = "$/caja 70 unidades"
$ws.Range("O171").Value = "Limache"
$ws.Range("P171").Value = 164
$ws.Range("Q171").Value = 70
$ws.Range("R171").Value = "Hortaliza"
